$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Bring the "Repayment schedule" sheet to the front (it becomes the active /
# selected tab after this edit).
$ws.Activate()

# Insert a new (blank) column before column N - this is the "Variable
# Instalments" column being added to the repayment schedule. All the old
# columns N:P (Late / Outstanding heading / Outstanding) shift one column
# to the right, becoming O:Q.
$ws.Columns("N").Insert()

# Match the column width used for the freshly inserted column.
$ws.Range("N1").ColumnWidth = 10.166666666666666

# Leave the cursor where the author left it when they made the edit.
$ws.Range("R7").Select() | Out-Null
